# Weekly cryptos-list refresh (GitHub Actions bot) — updates Price (D) and
# Volume(1h) (E) for each coin row, plus the three coin pairs that swapped
# rank order (rows 14/15 Solana<->BinanceUSD, 33/34 Filecoin<->Hedera,
# 40/41 TheSandbox<->Aptos) which also carry new Coin (B) / Link (C) text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Source cells are plain text (inlineStr) even when the text looks like a
# number ("0.9996", "317.14", ...). Assigning such a string straight to
# .Value lets Excel auto-coerce it to a Double, which would corrupt values
# like "28.443.08" style prices or drop trailing zeros ("1.112" -> 1.112 is
# fine, but "0.9996" -> 0.9996000000000000x is not). Force text by flipping
# NumberFormat to "@" for the write, then clear the format again so the
# cell keeps the workbook's original un-styled look.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $looksNumeric = $value -match '^-?\d+(\.\d+)?$'
    if ($looksNumeric) {
        $rng.NumberFormat = "@"
        $rng.Value = $value
        $rng.ClearFormats()
    } else {
        $rng.Value = $value
    }
}

$updates = [ordered]@{
    "D2" = '28.469.68'
    "E2" = '  +1.64%  '
    "D3" = '1.829.30'
    "E3" = '  +2.67%  '
    "D4" = '0.9996'
    "E4" = '  -0.05%  '
    "D5" = '317.14'
    "E5" = '  +0.31%  '
    "D6" = '0.9999'
    "E6" = '  +0.00%  '
    "D7" = '0.5346'
    "E7" = '  -0.84%  '
    "D8" = '0.4020'
    "E8" = '  +6.62%  '
    "D9" = '0.07655'
    "E9" = '  +2.82%  '
    "D10" = '41.83'
    "E10" = '  +0.46%  '
    "D11" = '1.112'
    "E11" = '  +1.74%  '
    "D12" = '6.330'
    "E12" = '  +4.17%  '
    "D13" = '7.622'
    "E13" = '  +5.52%  '
    "B14" = 'Solana'
    "C14" = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
    "D14" = '20.96'
    "E14" = '  +2.38%  '
    "B15" = 'BinanceUSD'
    "C15" = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
    "D15" = '0.9994'
    "E15" = '  -0.03%  '
    "D16" = '1.829.38'
    "E16" = '  +3.01%  '
    "D17" = '89.65'
    "E17" = '  +1.54%  '
    "D18" = '0.00001076'
    "E18" = '  +2.22%  '
    "D19" = '0.06588'
    "E19" = '  +2.33%  '
    "D20" = '17.72'
    "E20" = '  +2.95%  '
    "D22" = '6.080'
    "E22" = '  +3.49%  '
    "D23" = '28.458.44'
    "D24" = '11.19'
    "E24" = '  -0.26%  '
    "D25" = '2.218'
    "E25" = '  +6.32%  '
    "D26" = '2.465'
    "E26" = '  +8.30%  '
    "D27" = '157.23'
    "E27" = '  +0.81%  '
    "D28" = '20.67'
    "E28" = '  +2.28%  '
    "D29" = '2.039.84'
    "E29" = '  +3.06%  '
    "D30" = '124.14'
    "D31" = '1.125'
    "E31" = '  +1.80%  '
    "D32" = '0.1103'
    "E32" = '  +4.62%  '
    "B33" = 'Hedera'
    "C33" = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
    "D33" = '0.07542'
    "E33" = '  +17.12%  '
    "B34" = 'Filecoin'
    "C34" = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    "D34" = '5.676'
    "E34" = '  +2.92%  '
    "D35" = '3.644'
    "E35" = '  +0.07%  '
    "D36" = '0.2229'
    "E36" = '  -1.04%  '
    "D37" = '0.02341'
    "E37" = '  +2.57%  '
    "D38" = '5.235'
    "E38" = '  +4.64%  '
    "D39" = '8.838'
    "E39" = '  +4.93%  '
    "B40" = 'Aptos'
    "C40" = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
    "D40" = '11.31'
    "E40" = '  +2.16%  '
    "B41" = 'TheSandbox'
    "C41" = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
    "D41" = '0.6259'
    "E41" = '  +1.97%  '
    "D42" = '1.178'
    "E42" = '  +0.05%  '
    "D43" = '0.9997'
    "E43" = '  +0.04%  '
    "E44" = '  -3.61%  '
    "D45" = '13.48'
    "E45" = '  +2.12%  '
    "D46" = '3.699'
    "E46" = '  +0.80%  '
    "D47" = '0.5850'
    "E47" = '  +1.90%  '
    "D48" = '124.96'
    "E48" = '  -1.05%  '
    "D49" = '2.005'
    "E49" = '  +4.35%  '
    "D50" = '1.202'
    "E50" = '  +1.44%  '
    "D51" = '0.06893'
    "E51" = '  +1.49%  '
}

foreach ($cell in $updates.Keys) {
    Set-TextValue $cell $updates[$cell]
}
